# Remove the test/initialization card row (row 3, "test") from the sheet.
# This matches the commit message "Card __innit__ change to excel initialization":
# the placeholder "test" row used during initialization is deleted, shifting
# the remaining rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select the entire row before deleting it (mirrors the resulting selection
# state of A3:XFD3 once the row above shifts into row 3).
$ws.Rows(3).Select()
$ws.Rows(3).Delete()

# Leave the active cell/selection on row 3 (now occupied by the former row 4
# data), matching the saved selection of the whole row.
$ws.Range("A3:XFD3").Select()
